# Update building block types / metadata for ArrayExpress - Sequencing library template
$wb = $excel.ActiveWorkbook

# --- isa_template sheet: bump version number ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.3"

# --- Library sheet: update header names and term values ---
$wsLibrary = $wb.Worksheets.Item("Library")

# Header row (row 1)
$wsLibrary.Range("I1").Value = "Characteristic [library source]"
$wsLibrary.Range("R1").Value = "Output [Data]"

# Data row (row 2) - update ontology term URIs
$wsLibrary.Range("D2").Value = "https://bioregistry.io/EFO:0004184"
$wsLibrary.Range("H2").Value = "http://purl.org/nfdi4plants/ontology/dpbo/DPBO_0000086"
$wsLibrary.Range("K2").Value = "https://bioregistry.io/NCIT:C16629"
$wsLibrary.Range("N2").Value = "https://bioregistry.io/NCIT:C101294"
$wsLibrary.Range("Q2").Value = "https://bioregistry.io/NCIT:C17003"

# --- Update the table column names (annotationTable) to keep them in sync with headers ---
$table = $wsLibrary.ListObjects.Item("annotationTable")
$table.ListColumns.Item("Parameter [library source]").Name = "Characteristic [library source]"
$table.ListColumns.Item("Output [Raw Data File]").Name = "Output [Data]"

$wb.Save()
